$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 213, shifting the rest of the
# table (previously rows 213:336) down to 215:338.
$ws.Rows.Item(213).Insert()
$ws.Rows.Item(213).Insert()

$newDate = Get-Date -Year 2021 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0

# New row 213 ("Primera")
$ws.Cells.Item(213, 1).Value = 3
$ws.Cells.Item(213, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = $newDate
$ws.Cells.Item(213, 5).Value = 5
$ws.Cells.Item(213, 6).Value = 100112037
$ws.Cells.Item(213, 7).Value = "Cebollín"
$ws.Cells.Item(213, 8).Value = "Sin especificar"
$ws.Cells.Item(213, 9).Value = "Primera"
$ws.Cells.Item(213, 10).Value = 180
$ws.Cells.Item(213, 11).Value = 3000
$ws.Cells.Item(213, 12).Value = 3000
$ws.Cells.Item(213, 13).Value = 3000
$ws.Cells.Item(213, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(213, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(213, 16).Value = 83
$ws.Cells.Item(213, 17).Value = 36
$ws.Cells.Item(213, 18).Value = "Hortaliza"

# New row 214 ("Segunda")
$ws.Cells.Item(214, 1).Value = 3
$ws.Cells.Item(214, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(214, 3).Value = "Coquimbo"
$ws.Cells.Item(214, 4).Value = $newDate
$ws.Cells.Item(214, 5).Value = 5
$ws.Cells.Item(214, 6).Value = 100112037
$ws.Cells.Item(214, 7).Value = "Cebollín"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Segunda"
$ws.Cells.Item(214, 10).Value = 160
$ws.Cells.Item(214, 11).Value = 2000
$ws.Cells.Item(214, 12).Value = 2000
$ws.Cells.Item(214, 13).Value = 2000
$ws.Cells.Item(214, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(214, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(214, 16).Value = 56
$ws.Cells.Item(214, 17).Value = 36
$ws.Cells.Item(214, 18).Value = "Hortaliza"
